$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Wipe the old header row so no stray data survives (old sheet used A1:H1). ---
$ws.Range("A1:H1").ClearContents()

# --- New header row (A:G). ---
$ws.Range("A1").Value = "Month"
$ws.Range("B1").Value = "Remin Date"
$ws.Range("C1").Value = "AP NO."
$ws.Range("D1").Value = "Vender Name"
$ws.Range("E1").Value = "Transaction date"
$ws.Range("F1").Value = "Local Payday"
$ws.Range("G1").Value = "Amount"

# Center-align just the AP NO. header cell.
$ws.Range("C1").HorizontalAlignment = -4108

# --- Data rows 2-6. Columns B (Remin Date) and F (Local Payday) hold
#     numeric-looking date codes that must stay text, like the source file,
#     so force a Text number format right before the write and then drop
#     the format again (keeps the value/type without leaving a visible
#     style behind). ---
$textCells = @("B2","B3","B4","B5","B6","F2","F3","F4","F5","F6")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("B2").Value = "20201015"
$ws.Range("C2").Value = "20201005-ZSAC-0005"
$ws.Range("D2").Value = "AERO PERFORMANCE"
$ws.Range("F2").Value = "20201007"
$ws.Range("G2").Value = 55827.44

$ws.Range("B3").Value = "20201015"
$ws.Range("C3").Value = "20201005-ZSAC-0006"
$ws.Range("D3").Value = "DUNCAN AVIATION"
$ws.Range("F3").Value = "20201007"
$ws.Range("G3").Value = 1200

$ws.Range("B4").Value = "20201015"
$ws.Range("C4").Value = "20201005-ZSAC-0007"
$ws.Range("D4").Value = "WILLIS TOWERS WATSON INSURANCE SERVICES WEST, "
$ws.Range("F4").Value = "20201007"
$ws.Range("G4").Value = 26.83

$ws.Range("B5").Value = "20201015"
$ws.Range("C5").Value = "20201005-ZSAC-0008"
$ws.Range("D5").Value = "AIRCRAFT SPRUCE & SPECIALTY CO."
$ws.Range("F5").Value = "20201013"
$ws.Range("G5").Value = 1069.95

$ws.Range("B6").Value = "20201015"
$ws.Range("C6").Value = "20201005-ZSAC-0009"
$ws.Range("D6").Value = "ADAM JEHN DUNG TSUEI"
$ws.Range("F6").Value = "20201013"
$ws.Range("G6").Value = 3921.06

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

# --- Column widths, approximating the author's auto-fit pass. ---
$ws.Columns("B").ColumnWidth = 9.86
$ws.Columns("C").ColumnWidth = 18.29
$ws.Columns("D").ColumnWidth = 49.02
$ws.Columns("E").ColumnWidth = 14.57
$ws.Columns("F").ColumnWidth = 11.02

# --- View: zoom 100%, selection parked on D16 (matches the saved view). ---
$excel.ActiveWindow.Zoom = 100
$ws.Range("D16").Select()
